$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "slavery" now pairs with new lexical item "RESTRICT" instead of "LOSE" ---
# Copy the formatting already used for the "GAIN"-style entries (style index 3)
# from B15 onto B13, then overwrite its text.
$ws.Cells.Item(15, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4122)
$ws.Cells.Item(13, 2).Value = "RESTRICT"

# --- New lexicon entry: "economic" (added to BRI content lexicon) ---
# Row 16: economic / LOSE / -1 / TRUE
$ws.Cells.Item(16, 1).Value = "economic"
$ws.Cells.Item(16, 2).Value = "LOSE"
$ws.Cells.Item(16, 3).Value = -1
$ws.Cells.Item(16, 4).Value = $true

# Row 17: economic / GAIN / -1 / TRUE (styled like the other "GAIN" rows)
$ws.Cells.Item(15, 2).Copy()
$ws.Cells.Item(17, 2).PasteSpecial(-4122)
$ws.Cells.Item(17, 1).Value = "economic"
$ws.Cells.Item(17, 2).Value = "GAIN"
$ws.Cells.Item(17, 3).Value = -1
$ws.Cells.Item(17, 4).Value = $true
